$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-parsed as numbers by Excel, so they stay as literal text like the source data.
$ws.Range("D4:D5").NumberFormat = "@"
$ws.Range("D7:D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15:D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22:D23").NumberFormat = "@"
$ws.Range("D25:D43").NumberFormat = "@"
$ws.Range("D45:D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "27.459.99"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "1.832.94"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -3.09%  "
$ws.Range("D5").Value = "316.00"
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("D7").Value = "0.4307"
$ws.Range("E7").Value = "  -2.27%  "
$ws.Range("D8").Value = "0.3708"
$ws.Range("E8").Value = "  -2.88%  "
$ws.Range("D9").Value = "0.07278"
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("D10").Value = "0.8696"
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("D12").Value = "1.833.06"
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("D13").Value = "6.696"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").Value = "0.07093"
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("D16").Value = "88.15"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("E17").Value = "  -3.01%  "
$ws.Range("D18").Value = "0.000008937"
$ws.Range("E18").Value = "  -1.83%  "
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("D20").Value = "15.30"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").Value = "27.476.00"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("D22").Value = "5.183"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").Value = "10.92"
$ws.Range("E23").Value = "  -3.20%  "
$ws.Range("D24").Value = "2.059.46"
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("D25").Value = "2.007"
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("D26").Value = "153.77"
$ws.Range("E26").Value = "  -3.28%  "
$ws.Range("D27").Value = "18.52"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").Value = "2.166"
$ws.Range("E28").Value = "  +7.89%  "
$ws.Range("D29").Value = "5.308"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").Value = "117.65"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").Value = "0.08876"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("D32").Value = "1.209"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "0.7701"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").Value = "4.507"
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("D35").Value = "2.908"
$ws.Range("E35").Value = "  -3.83%  "
$ws.Range("D36").Value = "1.004"
$ws.Range("E36").Value = "  -2.91%  "
$ws.Range("D37").Value = "1.125"
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("D38").Value = "0.01969"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").Value = "0.05301"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").Value = "7.200"
$ws.Range("E40").Value = "  +3.92%  "
$ws.Range("D41").Value = "2.871"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.1680"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.5098"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").Value = "10.67"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "106.59"
$ws.Range("E46").Value = "  -3.81%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.4744"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").Value = "0.06429"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("E49").Value = "  -2.88%  "
$ws.Range("E50").Value = "  -2.48%  "
$ws.Range("D51").Value = "1.832"
$ws.Range("E51").Value = "  -2.99%  "
